$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035299986006433
$ws.Range("D2").Value = 1.030431791681564
$ws.Range("E2").Value = 1.043346011103842
$ws.Range("F2").Value = 1.051940005133803
$ws.Range("I2").Value = 1.032742580201409
$ws.Range("J2").Value = 1.040414515650139
$ws.Range("K2").Value = 1.033242821553265
$ws.Range("L2").Value = 1.046120136598059
$ws.Range("M2").Value = 1.054690116457603
$ws.Range("N2").Value = 1.017209706397926
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037085572483742
$ws.Range("D3").Value = 1.030900693442999
$ws.Range("E3").Value = 1.044933390792502
$ws.Range("F3").Value = 1.053623950861821
$ws.Range("I3").Value = 1.032875542783676
$ws.Range("J3").Value = 1.041839375531064
$ws.Range("K3").Value = 1.033520989168909
$ws.Range("L3").Value = 1.047516503596506
$ws.Range("M3").Value = 1.056184571918853
$ws.Range("N3").Value = 1.017699242577486
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038237622969961
$ws.Range("D4").Value = 1.03120157438988
$ws.Range("E4").Value = 1.045957297606657
$ws.Range("F4").Value = 1.054709086658705
$ws.Range("I4").Value = 1.032958114843868
$ws.Range("J4").Value = 1.042757787388789
$ws.Range("K4").Value = 1.033697734688799
$ws.Range("L4").Value = 1.048416339404178
$ws.Range("M4").Value = 1.05714664086608
$ws.Range("N4").Value = 1.018014344995279
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038721160813734
$ws.Range("D5").Value = 1.031327458059373
$ws.Range("E5").Value = 1.046386986916477
$ws.Range("F5").Value = 1.055164216449505
$ws.Range("I5").Value = 1.032991998688921
$ws.Range("J5").Value = 1.043143046334291
$ws.Range("K5").Value = 1.03377126039345
$ws.Range("L5").Value = 1.04879375497511
$ws.Range("M5").Value = 1.057549923421826
$ws.Range("N5").Value = 1.018146420970302
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038802303428561
$ws.Range("D6").Value = 1.031348558873825
$ws.Range("E6").Value = 1.04645908938786
$ws.Range("F6").Value = 1.055240572898708
$ws.Range("I6").Value = 1.032997639287794
$ws.Range("J6").Value = 1.043207684033205
$ws.Range("K6").Value = 1.033783560033608
$ws.Range("L6").Value = 1.048857073723594
$ws.Range("M6").Value = 1.057617568009553
$ws.Range("N6").Value = 1.018168574183351
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038244087085434
$ws.Range("D7").Value = 1.031203258838563
$ws.Range("E7").Value = 1.045963042108931
$ws.Range("F7").Value = 1.054715172276793
$ws.Range("I7").Value = 1.03295857086078
$ws.Range("J7").Value = 1.042762938526371
$ws.Range("K7").Value = 1.033698720200941
$ws.Range("L7").Value = 1.048421385872697
$ws.Range("M7").Value = 1.057152034132527
$ws.Range("N7").Value = 1.018016111339447
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035904135403513
$ws.Range("D8").Value = 1.030590782454567
$ws.Range("E8").Value = 1.043883151059292
$ws.Range("F8").Value = 1.052510040305762
$ws.Range("I8").Value = 1.032788232562463
$ws.Range("J8").Value = 1.040896800155646
$ws.Range("K8").Value = 1.033337501658291
$ws.Range("L8").Value = 1.046592820231973
$ws.Range("M8").Value = 1.055196205751611
$ws.Range("N8").Value = 1.017375493855205
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03175450176062
$ws.Range("D9").Value = 1.029492186141551
$ws.Range("E9").Value = 1.040192755209057
$ws.Range("F9").Value = 1.048589315768057
$ws.Range("I9").Value = 1.03246155499191
$ws.Range("J9").Value = 1.037580524830747
$ws.Range("K9").Value = 1.032676140422534
$ws.Range("L9").Value = 1.043341716069383
$ws.Range("M9").Value = 1.051711365178595
$ws.Range("N9").Value = 1.01623373022904
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028969336290308
$ws.Range("D10").Value = 1.028746836428481
$ws.Range("E10").Value = 1.037714601143161
$ws.Range("F10").Value = 1.045951097052791
$ws.Range("I10").Value = 1.03222594312556
$ws.Range("J10").Value = 1.03535010217368
$ws.Range("K10").Value = 1.032218551451625
$ws.Range("L10").Value = 1.04115408481449
$ws.Range("M10").Value = 1.049361514270864
$ws.Range("N10").Value = 1.015463591581948
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027758651087324
$ws.Range("D11").Value = 1.028421027322378
$ws.Range("E11").Value = 1.036637100945013
$ws.Range("F11").Value = 1.044802734588821
$ws.Range("I11").Value = 1.032119691045545
$ws.Range("J11").Value = 1.034379479095911
$ws.Range("K11").Value = 1.032016456680416
$ws.Range("L11").Value = 1.040201844266875
$ws.Range("M11").Value = 1.048337508793635
$ws.Range("N11").Value = 1.015127923923465
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027308223880283
$ws.Range("D12").Value = 1.028299546813682
$ws.Range("E12").Value = 1.036236185412986
$ws.Range("F12").Value = 1.044375263549758
$ws.Range("I12").Value = 1.032079588386593
$ws.Range("J12").Value = 1.034018204061131
$ws.Range("K12").Value = 1.031940795391453
$ws.Range("L12").Value = 1.039847376081621
$ws.Range("M12").Value = 1.04795615462417
$ws.Range("N12").Value = 1.015002907081969
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027404875171458
$ws.Range("D13").Value = 1.02832562562831
$ws.Range("E13").Value = 1.036322214373855
$ws.Range("F13").Value = 1.044466999363247
$ws.Range("I13").Value = 1.032088219324764
$ws.Range("J13").Value = 1.034095732636423
$ws.Range("K13").Value = 1.031957051888458
$ws.Range("L13").Value = 1.039923445513542
$ws.Range("M13").Value = 1.048038001575628
$ws.Range("N13").Value = 1.01502973886455
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027721433537619
$ws.Range("D14").Value = 1.028410995097715
$ws.Range("E14").Value = 1.036603975191155
$ws.Range("F14").Value = 1.044767418528643
$ws.Range("I14").Value = 1.032116389121548
$ws.Range("J14").Value = 1.034349631226814
$ws.Range("K14").Value = 1.032010214618294
$ws.Range("L14").Value = 1.040172559482984
$ws.Range("M14").Value = 1.048306006286542
$ws.Range("N14").Value = 1.015117596857012
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027916378941239
$ws.Range("D15").Value = 1.028463533035007
$ws.Range("E15").Value = 1.036777486199293
$ws.Range("F15").Value = 1.044952394516385
$ws.Range("I15").Value = 1.03213366118906
$ws.Range("J15").Value = 1.034505967670824
$ws.Range("K15").Value = 1.032042891180095
$ws.Range("L15").Value = 1.040325945172283
$ws.Range("M15").Value = 1.048471000785417
$ws.Range("N15").Value = 1.015171684518071
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029049584947948
$ws.Range("D16").Value = 1.028768394707432
$ws.Range("E16").Value = 1.03778601627211
$ws.Range("F16").Value = 1.0460271822889
$ws.Range("I16").Value = 1.032232905567818
$ws.Range("J16").Value = 1.035414415999239
$ws.Range("K16").Value = 1.032231880476163
$ws.Range("L16").Value = 1.041217175631637
$ws.Range("M16").Value = 1.049429335763923
$ws.Range("N16").Value = 1.01548582207505
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029759145290885
$ws.Range("D17").Value = 1.028958805360122
$ws.Range("E17").Value = 1.038417439449044
$ws.Range("F17").Value = 1.046699750808507
$ws.Range("I17").Value = 1.032294026099777
$ws.Range("J17").Value = 1.035982955886292
$ws.Range("K17").Value = 1.03234936944118
$ws.Range("L17").Value = 1.041774876535409
$ws.Range("M17").Value = 1.050028721520119
$ws.Range("N17").Value = 1.015682281385862
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030172568676995
$ws.Range("D18").Value = 1.029069572657465
$ws.Range("E18").Value = 1.0387853103313
$ws.Range("F18").Value = 1.047091471263125
$ws.Range("I18").Value = 1.032329268507136
$ws.Range("J18").Value = 1.03631411013192
$ws.Range("K18").Value = 1.032417517101583
$ws.Range("L18").Value = 1.042099694405879
$ws.Range("M18").Value = 1.050377706573569
$ws.Range("N18").Value = 1.015796661606911
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030313459449814
$ws.Range("D19").Value = 1.029107291220223
$ws.Range("E19").Value = 1.038910672765513
$ws.Range("F19").Value = 1.047224940396985
$ws.Range("I19").Value = 1.032341216033393
$ws.Range("J19").Value = 1.036426946702209
$ws.Range("K19").Value = 1.032440688953612
$ws.Range("L19").Value = 1.042210368007932
$ws.Range("M19").Value = 1.050496595661091
$ws.Range("N19").Value = 1.015835626653956
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029683062991994
$ws.Range("D20").Value = 1.028938406714372
$ws.Range("E20").Value = 1.038349738062469
$ws.Range("F20").Value = 1.046627650416998
$ws.Range("I20").Value = 1.032287510665678
$ws.Range("J20").Value = 1.035922005195823
$ws.Range("K20").Value = 1.032336803464371
$ws.Range("L20").Value = 1.041715090226457
$ws.Range("M20").Value = 1.049964478017836
$ws.Range("N20").Value = 1.015661225043316
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027628235160652
$ws.Range("D21").Value = 1.02838586862605
$ws.Range("E21").Value = 1.036521022667742
$ws.Range("F21").Value = 1.044678978055448
$ws.Range("I21").Value = 1.032108111375411
$ws.Range("J21").Value = 1.034274885043486
$ws.Range("K21").Value = 1.031994575925108
$ws.Range("L21").Value = 1.040099222831505
$ws.Range("M21").Value = 1.048227113116484
$ws.Range("N21").Value = 1.01509173415843
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026332079586497
$ws.Range("D22").Value = 1.028035801822122
$ws.Range("E22").Value = 1.035367271080526
$ws.Range("F22").Value = 1.043448449898513
$ws.Range("I22").Value = 1.03199163656742
$ws.Range("J22").Value = 1.033234971825414
$ws.Range("K22").Value = 1.03177596543996
$ws.Range("L22").Value = 1.039078836949166
$ws.Range("M22").Value = 1.047129010026396
$ws.Range("N22").Value = 1.014731732053806
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027019601454832
$ws.Range("D23").Value = 1.028221631281289
$ws.Range("E23").Value = 1.035979278078013
$ws.Range("F23").Value = 1.044101286293346
$ws.Range("I23").Value = 1.032053730946782
$ws.Range("J23").Value = 1.033786662795887
$ws.Range("K23").Value = 1.03189218093694
$ws.Range("L23").Value = 1.039620187467992
$ws.Range("M23").Value = 1.047711685863369
$ws.Range("N23").Value = 1.014922761790638
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029717442701623
$ws.Range("D24").Value = 1.028947624899415
$ws.Range("E24").Value = 1.038380330725869
$ws.Range("F24").Value = 1.046660231264331
$ws.Range("I24").Value = 1.032290455971302
$ws.Range("J24").Value = 1.035949547627232
$ws.Range("K24").Value = 1.03234248266802
$ws.Range("L24").Value = 1.04174210656784
$ws.Range("M24").Value = 1.049993508826968
$ws.Range("N24").Value = 1.015670740150592
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032830508487089
$ws.Range("D25").Value = 1.029778485672828
$ws.Range("E25").Value = 1.041149902068874
$ws.Range("F25").Value = 1.049607153469448
$ws.Range("I25").Value = 1.03254915108397
$ws.Range("J25").Value = 1.038441250930218
$ws.Range("K25").Value = 1.032850059893946
$ws.Range("L25").Value = 1.044185712507498
$ws.Range("M25").Value = 1.052616913775952
$ws.Range("N25").Value = 1.0165304624545
